$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.061.90"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").Value = "1.601.90"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "212.33"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "18.09"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  +4.55%  "
$ws.Range("D12").Value = "1.825.76"
$ws.Range("E12").Value = "  +3.13%  "
$ws.Range("D13").Value = "1.599.48"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "26.063.10"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").Value = "60.36"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "201.13"
$ws.Range("E20").Value = "  +8.51%  "
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "5.99"
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("E24").Value = "  +12.69%  "
$ws.Range("D25").Value = "141.30"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  -5.17%  "
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("D33").Value = "2.96"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").Value = "1.124.03"
$ws.Range("E36").Value = "  +3.85%  "
$ws.Range("E37").Value = "  +10.15%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "0.790"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "1.737.91"
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("D45").Value = "93.16"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  +3.29%  "
$ws.Range("D47").Value = "53.34"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "0.408"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "0.0₇0923"
$ws.Range("E51").Value = "  -16.61%  "
